$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.51
$ws.Range("G2").Value = 1.74
$ws.Range("H2").Value = 2.34
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 3.6
$ws.Range("K2").Value = 6
$ws.Range("P2").Value = 1.86
$ws.Range("Q2").Value = 1.9

# Row 3
$ws.Range("F3").Value = 1.74
$ws.Range("G3").Value = 1.93
$ws.Range("P3").Value = 1.71

# Row 4
$ws.Range("G4").Value = 2.02
$ws.Range("H4").Value = 3.75
$ws.Range("J4").Value = 3.1
$ws.Range("P4").Value = 1.6
$ws.Range("Q4").Value = 2.28
